$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 72) to the bottom of the portfolio table,
# carrying forward the prior row's quote values under the next date.
$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "2025-10-26"
$ws.Range("A72").Style = "Normal"

$ws.Range("B72").Value = 53.81999969482422
$ws.Range("C72").Value = 403.2999877929688
$ws.Range("D72").Value = 326.6000061035156
